$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.590.35'
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").Value = '1.742.71'
$ws.Range("E3").Value = '  +0.82%  '

$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '''246.59'
$ws.Range("E5").Value = '  +0.86%  '

$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("D7").Value = '''0.4924'
$ws.Range("E7").Value = '  +2.26%  '

$ws.Range("D8").Value = '''0.2673'
$ws.Range("E8").Value = '  -0.43%  '

$ws.Range("D9").Value = '''0.06291'
$ws.Range("E9").Value = '  +0.95%  '

$ws.Range("D10").Value = '1.742.61'
$ws.Range("E10").Value = '  +0.78%  '

$ws.Range("D11").Value = '''0.07051'
$ws.Range("E11").Value = '  -1.07%  '

$ws.Range("D12").Value = '''15.74'
$ws.Range("E12").Value = '  -0.08%  '

$ws.Range("D13").Value = '''0.6155'
$ws.Range("E13").Value = '  -0.52%  '

$ws.Range("D14").Value = '''4.588'
$ws.Range("E14").Value = '  +0.69%  '

$ws.Range("E15").Value = '  +1.31%  '

$ws.Range("E16").Value = '  +0.02%  '

$ws.Range("D17").Value = '26.598.39'
$ws.Range("E17").Value = '  +0.46%  '

$ws.Range("B18").Value = 'BinanceUSD'
$ws.Range("C18").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D18").Value = '''1.001'
$ws.Range("E18").Value = '  +0.10%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '''0.000007292'
$ws.Range("E19").Value = '  +4.81%  '

$ws.Range("D20").Value = '''11.56'
$ws.Range("E20").Value = '  -1.43%  '

$ws.Range("D21").Value = '1.966.12'
$ws.Range("E21").Value = '  +0.56%  '

$ws.Range("D22").Value = '''4.581'
$ws.Range("E22").Value = '  +0.55%  '

$ws.Range("D23").Value = '''8.727'
$ws.Range("E23").Value = '  -2.28%  '

$ws.Range("D24").Value = '''5.268'
$ws.Range("E24").Value = '  -0.97%  '

$ws.Range("D25").Value = '''139.38'
$ws.Range("E25").Value = '  +2.11%  '

$ws.Range("D26").Value = '''15.52'
$ws.Range("E26").Value = '  +0.96%  '

$ws.Range("D27").Value = '''1.426'
$ws.Range("E27").Value = '  +1.46%  '

$ws.Range("D28").Value = '''1.764'
$ws.Range("E28").Value = '  -1.97%  '

$ws.Range("D29").Value = '''107.70'
$ws.Range("E29").Value = '  +0.85%  '

$ws.Range("D30").Value = '''4.049'
$ws.Range("E30").Value = '  +1.49%  '

$ws.Range("D31").Value = '''0.08052'
$ws.Range("E31").Value = '  +0.35%  '

$ws.Range("D32").Value = '''3.734'
$ws.Range("E32").Value = '  -0.43%  '

$ws.Range("D33").Value = '''0.04633'
$ws.Range("E33").Value = '  +1.59%  '

$ws.Range("E34").Value = '  -0.20%  '

$ws.Range("D35").Value = '''1.017'
$ws.Range("E35").Value = '  +2.64%  '

$ws.Range("D36").Value = '''0.6375'
$ws.Range("E36").Value = '  -0.64%  '

$ws.Range("E37").Value = '  +3.22%  '

$ws.Range("D38").Value = '''0.8991'
$ws.Range("E38").Value = '  -4.58%  '

$ws.Range("D39").Value = '''2.422'
$ws.Range("E39").Value = '  +0.39%  '

$ws.Range("D40").Value = '''1.003'
$ws.Range("E40").Value = '  +0.05%  '

$ws.Range("D41").Value = '''0.01505'
$ws.Range("E41").Value = '  +0.20%  '

$ws.Range("D42").Value = '''101.84'
$ws.Range("E42").Value = '  -4.73%  '

$ws.Range("D43").Value = '''5.422'
$ws.Range("E43").Value = '  -4.34%  '

$ws.Range("D44").Value = '''0.3921'
$ws.Range("E44").Value = '  -0.02%  '

$ws.Range("D45").Value = '''6.886'
$ws.Range("E45").Value = '  -2.02%  '

$ws.Range("D46").Value = '''0.1183'
$ws.Range("E46").Value = '  -1.11%  '

$ws.Range("D47").Value = '''0.05402'
$ws.Range("E47").Value = '  +1.51%  '

$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = '''30.57'
$ws.Range("E48").Value = '  -1.37%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''7.810'
$ws.Range("E49").Value = '  -0.74%  '

$ws.Range("D50").Value = '''1.264'
$ws.Range("E50").Value = '  -0.57%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '''51.69'
$ws.Range("E51").Value = '  +0.57%  '
